# Applies the "Added more techs to network, get csp fraction for plt" edit
# to CSP_case.xlsx (sheet "case_input_intermodel_one_node_").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) CASE_DATA block updates
# ---------------------------------------------------------------------
$ws.Range("B26").Value = "input_files/costs_concentrated_solar.csv"   # costs_path
$ws.Range("B28").Value = "csp_test"                                   # case_name
$ws.Range("B29").Value = "csp_test_case"                              # filename_prefix
$ws.Range("B30").Value = "2023-01-01 00:00:00"                        # datetime_start (unchanged)
$ws.Range("B31").Value = "2023-12-31 23:00:00"                        # datetime_end (extended to full year)

# ---------------------------------------------------------------------
# 2) Make room in COMPONENT_DATA for the new technologies: insert 6
#    fresh rows right after the existing csp Generator row (51) and
#    before the old "lost_load" Generator row (so everything from the
#    old row 52 onward - lost_load, END_COMPONENT_DATA, the trailing
#    notes - shifts down by 6 rows, matching the target layout).
# ---------------------------------------------------------------------
$ws.Rows("52:57").Insert()

# ---------------------------------------------------------------------
# 3) Update the existing Load row (50): bus switches from electricity
#    to heat.
# ---------------------------------------------------------------------
$ws.Range("D50").Value = "heat"

# ---------------------------------------------------------------------
# 4) Update the existing csp Generator row (51): rename, change bus to
#    heat, drop p_max_pu (H51), and replace the capital/marginal cost /
#    efficiency / lifetime numbers with "db" placeholders (now sourced
#    from the external cost database instead of being hard-coded).
# ---------------------------------------------------------------------
$ws.Range("B51").Value = "csp glasspoint"
$ws.Range("C51").Value = "concentrated solar"
$ws.Range("D51").Value = "heat"
$ws.Range("H51").ClearContents()
$ws.Range("I51").Value = "db"
$ws.Range("K51").Value = "db"
$ws.Range("N51").Value = "db"
$ws.Range("O51").Value = "db"

# ---------------------------------------------------------------------
# 5) Populate the 6 freshly inserted rows (52-57) with the new
#    components: molten-salt storage charger/store/discharger, a gas
#    boiler, a second solar generator, and an electric boiler link.
# ---------------------------------------------------------------------

# Row 52: Molten-Salt charger (Link)
$ws.Range("A52").Value = "#Link"
$ws.Range("B52").Value = "Molten-Salt-charger glasspoint"
$ws.Range("C52").Value = "molten salt charger"
$ws.Range("D52").Value = "heat"
$ws.Range("E52").Value = "salt"
$ws.Range("I52").Value = "db"
$ws.Range("K52").Value = "db"
$ws.Range("N52").Value = "db"
$ws.Range("O52").Value = "db"

# Row 53: Molten-Salt store (Store)
$ws.Range("A53").Value = "#Store"
$ws.Range("B53").Value = "Molten-Salt-store glasspoint"
$ws.Range("C53").Value = "molten salt storage"
$ws.Range("D53").Value = "salt"
$ws.Range("I53").Value = "db"
$ws.Range("K53").Value = "db"
$ws.Range("M53").Value = $true
$ws.Range("N53").Value = "db"
$ws.Range("O53").Value = "db"

# Row 54: Molten-Salt discharger (Link)
$ws.Range("A54").Value = "#Link"
$ws.Range("B54").Value = "Molten-Salt-discharger glasspoint"
$ws.Range("C54").Value = "molten salt disbicharger"
$ws.Range("D54").Value = "heat"
$ws.Range("E54").Value = "salt"
$ws.Range("I54").Value = "db"
$ws.Range("K54").Value = "db"
$ws.Range("N54").Value = "db"
$ws.Range("O54").Value = "db"

# Row 55: Gas boiler steam (Generator)
$ws.Range("A55").Value = "Generator"
$ws.Range("B55").Value = "gas boiler steam"
$ws.Range("C55").Value = "gas boiler"
$ws.Range("D55").Value = "heat"
$ws.Range("I55").Value = "db"
$ws.Range("K55").Value = "db"
$ws.Range("N55").Value = "db"
$ws.Range("O55").Value = "db"

# Row 56: solar-utility single-axis tracking (Generator)
$ws.Range("A56").Value = "Generator"
$ws.Range("B56").Value = "solar-utility single-axis tracking"
$ws.Range("C56").Value = "solar"
$ws.Range("D56").Value = "electricity"
$ws.Range("I56").Value = "db"
$ws.Range("K56").Value = "db"
$ws.Range("N56").Value = "db"
$ws.Range("O56").Value = "db"

# Row 57: electric boiler steam (Link)
$ws.Range("A57").Value = "Link"
$ws.Range("B57").Value = "electric boiler steam"
$ws.Range("C57").Value = "electric boiler"
$ws.Range("D57").Value = "electricity"
$ws.Range("E57").Value = "heat"
$ws.Range("I57").Value = "db"
$ws.Range("K57").Value = "db"
$ws.Range("N57").Value = "db"
$ws.Range("O57").Value = "db"

# ---------------------------------------------------------------------
# 6) Row 58 is the old "lost_load" Generator row, shifted down by the
#    insert above; that technology no longer exists in the new case,
#    so clear it out - it becomes a blank template row like row 50/58
#    used to be.
# ---------------------------------------------------------------------
$ws.Range("A58").ClearContents()
$ws.Range("B58").ClearContents()
$ws.Range("D58").ClearContents()
$ws.Range("K58").ClearContents()
$ws.Range("L58").ClearContents()

# ---------------------------------------------------------------------
# 7) Restore the active selection to match the edited workbook.
# ---------------------------------------------------------------------
$ws.Range("B31").Select()
